# "Fruta / hortaliza, semanal" weekly update.
# A new weekly price record is inserted as row 107 (pushing the existing
# rows 107-132 down to 108-133), adding one more "Zapallo italiano" /
# "Vega Monumental Concepción" observation dated 2022-06-24 (serial 44736).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 107..132 down to 108..133, leaving a blank row 107.
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A107").Value = 11
$ws.Range("B107").Value = "Vega Monumental Concepción"
$ws.Range("C107").Value = "Bíobío"
$ws.Range("D107").Value = 44736
$ws.Range("E107").Value = 8
$ws.Range("F107").Value = 100112032
$ws.Range("G107").Value = "Zapallo italiano"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 16000
$ws.Range("L107").Value = 17000
$ws.Range("M107").Value = 16500
$ws.Range("N107").Value = "$/caja 60 unidades"
$ws.Range("O107").Value = "Región de Arica y Parinacota"
$ws.Range("P107").Value = 275
$ws.Range("Q107").Value = 60
$ws.Range("R107").Value = "Hortaliza"
